# Add a new "2022-Q4" sheet (feat: add 2022-Q4 data), placed right after
# "总计" and before "2022-Q3", and insert the corresponding summary row at
# the top of the "总计" (总计) sheet's data table.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Create the "2022-Q4" worksheet by cloning "2022-Q3" so it inherits
#    the exact same layout/styling (bold+bordered header row, bordered
#    index column), then insert it immediately before "2022-Q3".
# ------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3Index = $q3.Index
$q3.Copy($q3, $null)

# Copy() inserts the clone at the original sheet's index and pushes the
# source sheet (keeping its original name) one slot later, so the new
# clone - still carrying the auto "(2)" suffix - is the sheet now sitting
# at $q3Index.
$q4 = $wb.Worksheets.Item($q3Index)
$q4.Name = "2022-Q4"

# The cloned sheet has 10 data rows (rows 2-11); 2022-Q4 only has 7
# funds (rows 2-8), so drop the trailing 3 rows.
$q4.Rows.Item(9).Resize(3).Delete()

# ------------------------------------------------------------------
# 2) Populate the 2022-Q4 fund table.
#    Columns: A idx, B code, C name, D size, E stock-position,
#             F position-pct, G holding-value, H position-rank
#    B-G are stored as text in the source data (leading/trailing
#    zeros matter), so write them with a leading apostrophe to force
#    text and avoid Excel's automatic number coercion. H and A stay
#    numeric.
# ------------------------------------------------------------------
$q4Data = @(
    @(0, "501208", "中欧创新未来混合（LOF）", "54.76", "85.24", "3.60", "1.9714", 8),
    @(1, "005763", "中欧电子信息产业沪港深股票C", "14.81", "91.56", "5.79", "0.8575", 4),
    @(2, "004616", "中欧电子信息产业沪港深股票A", "6.80", "91.56", "5.79", "0.3937", 4),
    @(3, "011868", "中信建投远见回报混合A", "6.34", "95.01", "3.43", "0.2175", 10),
    @(4, "011869", "中信建投远见回报混合C", "1.41", "95.01", "3.43", "0.0484", 10),
    @(5, "015412", "西部利得数字产业混合A", "0.94", "92.30", "5.03", "0.0473", 8),
    @(6, "015413", "西部利得数字产业混合C", "0.43", "92.30", "5.03", "0.0216", 8)
)

$r = 2
foreach ($row in $q4Data) {
    $q4.Cells.Item($r, 1).Value = $row[0]
    $q4.Cells.Item($r, 2).Value = "'" + $row[1]
    $q4.Cells.Item($r, 3).Value = "'" + $row[2]
    $q4.Cells.Item($r, 4).Value = "'" + $row[3]
    $q4.Cells.Item($r, 5).Value = "'" + $row[4]
    $q4.Cells.Item($r, 6).Value = "'" + $row[5]
    $q4.Cells.Item($r, 7).Value = "'" + $row[6]
    $q4.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# ------------------------------------------------------------------
# 3) Insert a new top data row in "总计" for 2022-Q4 (10 持有数量, 亿元
#    持有市值 totals), pushing the rest of the quarters down by one row.
#    Shift existing rows 2..6 down to 3..7 by copying values (bottom to
#    top, to avoid clobbering), rather than a structural row Insert(),
#    which would drag unwanted formatting from neighbouring rows onto
#    the new cells.
# ------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

for ($i = 5; $i -ge 2; $i--) {
    $srcRow = $i
    $dstRow = $i + 1
    $summary.Cells.Item($dstRow, 2).Value = $summary.Cells.Item($srcRow, 2).Value
    $summary.Cells.Item($dstRow, 3).Value = $summary.Cells.Item($srcRow, 3).Value
    $summary.Cells.Item($dstRow, 4).Value = $summary.Cells.Item($srcRow, 4).Value
    $summary.Cells.Item($dstRow, 1).Value = $summary.Cells.Item($srcRow, 1).Value + 1
}

# Write the new 2022-Q4 row into row 2.
$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q4"
$summary.Cells.Item(2, 3).Value = 7
$summary.Cells.Item(2, 4).Value = 3.56

# Column A keeps the bold+bordered header style throughout; make sure
# the newly written A2 carries it too (it already should, since we
# never touched the cell's formatting - only its value).
